# ScrabbleBoardLayout.xlsx update
# - Narrow the board columns (A:O) from 4.6640625 to 3.109375
# - Place the letters of the finished demo words on the board
#   (write order matters: it controls the order strings are interned
#   into the shared-strings table, matching the target file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize columns A through O (target stored width 3.109375; the host's
# ColumnWidth->stored-width conversion snaps to the nearest 1/6, so 14/6
# is the closest input that reproduces it)
$ws.Range("A1:O1").EntireColumn.ColumnWidth = 2.3333333333333335

# Fill in the letters - order chosen to reproduce the shared-strings
# interning order of the target workbook (first-seen order of each
# unique letter, then the remaining duplicate cells).
$ws.Range("E15").Value = "D"
$ws.Range("G8").Value = "O"
$ws.Range("H10").Value = "S"
$ws.Range("F8").Value = "H"
$ws.Range("H8").Value = "R"
$ws.Range("I8").Value = "N"
$ws.Range("H6").Value = "F"
$ws.Range("H7").Value = "A"
$ws.Range("H9").Value = "M"
$ws.Range("F10").Value = "P"
$ws.Range("I10").Value = "T"
$ws.Range("J10").Value = "E"
$ws.Range("J9").Value = "B"
$ws.Range("F11").Value = "I"
$ws.Range("I9").Value = "O"
$ws.Range("G10").Value = "A"
$ws.Range("E11").Value = "B"
$ws.Range("G11").Value = "T"
$ws.Range("E12").Value = "O"
$ws.Range("E13").Value = "A"
$ws.Range("E14").Value = "R"
$ws.Range("A15").Value = "H"
$ws.Range("B15").Value = "E"
$ws.Range("C15").Value = "A"
$ws.Range("D15").Value = "R"
